$wb = $excel.ActiveWorkbook
$wsMain = $wb.Worksheets.Item("TestStepExecution")
$wsData = $wb.Worksheets.Item("data")

# Add 8 more "Dummy" entries to the data sheet (rows 4-11),
# matching the existing rows 1-3.
for ($i = 4; $i -le 11; $i++) {
    $wsData.Cells.Item($i, 1).Value = "Dummy"
}

# After typing, Excel's cursor rests one row below the last entry.
[void]$wsData.Cells.Item(12, 1).Select()

# Restore the originally active sheet/tab.
[void]$wsMain.Activate()
